## fix v-for của data.js 17:40 25/03/21
##
## 1. Sheet1 ("Màn hình chính"): insert two new test-case rows (row 6 and
##    row 7 in the final layout) right after the existing "Font chữ..."
##    row, pushing every row below down by one.
## 2. Sheet2 ("Màn hình thêm mới") and Sheet3 ("Màn hình sửa"): fill in
##    the previously-blank row 3 with a new test case, and shrink the
##    oversized row (sheet2 row 8 / sheet3 row 1).
## 3. Refresh each sheet's remembered selection, and make "Màn hình
##    chính" (sheet 1) the active tab again (it was "Popup xóa").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Màn hình chính"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Push rows 7.. down by one empty row, then add two real rows (6 and 7)
# of new content in the gap that opens up.
[void]$ws1.Rows.Item(6).Insert()

$ws1.Range("B6").Value = "Hiển thị chữ"
$ws1.Range("E6").Value = "Chữ căn trái"

$ws1.Range("B7").Value = "Hiển thị số"
$ws1.Range("E7").Value = "Số căn phải, giá tiền có chấm "
$ws1.Range("F7").Value = "Nhập tối đa 20 ký tự"

# ---------------------------------------------------------------------
# Sheet 2: "Màn hình thêm mới"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B3").Value = "Hiển thị giá tiền"
$ws2.Range("E3").Value = "Giá tiền có chấm, căn phải"
$ws2.Range("F3").Value = "Nhập tối đa được 20 ký tự"

# Row 8 was oversized (75pt); shrink it back down to 30pt.
$ws2.Rows.Item(8).RowHeight = 30

# ---------------------------------------------------------------------
# Sheet 3: "Màn hình sửa"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Header row was oversized (37.5pt); shrink it back down to 18.75pt.
$ws3.Rows.Item(1).RowHeight = 18.75

$ws3.Range("B3").Value = "Hiển thị giá tiền"
$ws3.Range("E3").Value = "Giá tiền có chấm, căn phải"
$ws3.Range("F3").Value = "Nhập tối đa được 20 ký tự"

# ---------------------------------------------------------------------
# Sheet 4: "Popup xóa" -- no content changes, selection stays at J8.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
[void]$ws4.Range("J8").Select()

# ---------------------------------------------------------------------
# Restore each sheet's own remembered selection ...
# ---------------------------------------------------------------------
[void]$ws2.Rows.Item(3).Select()
[void]$ws3.Range("G3").Select()

# ... and finish with "Màn hình chính" as the active sheet/selection,
# matching the workbook reopening on its first tab.
[void]$ws1.Range("I9").Select()
